$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text values for the "Terv" (plan) column, rows 5-15
$ws.Range("B5").Value = "Use Case Diagram, Wireframe, adatbázis drótháló"
$ws.Range("B6").Value = "Use Case Diagram, Wireframe, adatbázis drótháló"
$ws.Range("B7").Value = "Adatbázis implementáció"
$ws.Range("B8").Value = "Adatbázis implementáció"
$ws.Range("B9").Value = "Backend Implementáció"
$ws.Range("B10").Value = "Backend Implementáció"
$ws.Range("B11").Value = "Frontend implementáció"
$ws.Range("B12").Value = "Frontend implementáció"
$ws.Range("B13").Value = "Javítás, összeillesztési problémák feloldása"
$ws.Range("B14").Value = "Diagramok készítése, Dokumentáció, (0.Leadás) "
$ws.Range("B15").Value = "Visszajelzések alapján javítás, Hivatalos leadás"

# Apply "Wrap Text" to the whole sheet (all cells), matching the cols style="2" change
$ws.Cells.WrapText = $true

# Move the active selection to H13, matching the new selection in the diff
$ws.Range("H13").Select()
